# Update view-count column (F) for specific rows on the "展览" sheet
# and on the "全部类型" sheet (which aggregates the same records).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - row => new F value
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 1188
$ws1.Range("F8").Value  = 73
$ws1.Range("F12").Value = 274
$ws1.Range("F21").Value = 7305
$ws1.Range("F24").Value = 3318
$ws1.Range("F31").Value = 1341
$ws1.Range("F36").Value = 1418

# Sheet "全部类型" (All types) - same records appear at different rows
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F8").Value  = 1188
$ws4.Range("F11").Value = 73
$ws4.Range("F15").Value = 274
$ws4.Range("F25").Value = 7305
$ws4.Range("F28").Value = 3318
$ws4.Range("F36").Value = 1341
$ws4.Range("F41").Value = 1418
